# new layout G3 graphics
# Swap the "soft" (E) and "rigid" (F) columns: header labels and all data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the used range to know how many rows to process.
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $eCell = $ws.Cells.Item($r, 5)   # column E
    $fCell = $ws.Cells.Item($r, 6)   # column F

    $eVal = $eCell.Value2
    $fVal = $fCell.Value2

    if ($eVal -eq $null -and $fVal -eq $null) {
        continue
    }

    $eCell.Value2 = $fVal
    $fCell.Value2 = $eVal
}
